$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Update project title typo: फ्रेंश -> फ्रेंस
$ws.Range("A6").Value = "Project:- फ्रेंस कृषि तथा पशुपालन"

# 2. Clear the date text, keep the "Date:" label with trailing spaces
$ws.Range("H7").Value = "Date:                     "

# 3. Update CGI sheet thickness wording: 0.24mm (@$) -> 0.26mm (@^)
$ws.Range("B9").Value = ")=@^ dL=dL= afSnf] ;L=hL=cfO{= 5fgf 5fpg] sfd "

# 4. Update the width formula for the second flooring measurement: 15.33 -> 12
$ws.Range("D20").Formula = "=12/3.281"

# 5. Insert a new row at 29 (shifts old rows 29-33 down to 30-34),
#    duplicating the "Total Estimated" row as a new hidden helper row.
$ws.Rows.Item(29).Insert()
$ws.Range("B29").Value = "Total Estimated"
$ws.Range("C29").Formula = "=J26"
$ws.Range("E29").Value = 100
$ws.Rows.Item(29).Hidden = $true

# 6. Clear the old total row's percentage value (E28 no longer holds 100)
$ws.Range("E28").ClearContents()

# 7. Update the shifted percentage/remainder formulas to reference the new
#    "Total Estimated" helper row (C29) instead of the original row (C28).
$ws.Range("E31").Formula = "=C31/C29*100"
$ws.Range("C32").Formula = "=C29-C31"

# 8. Hide the budget-breakdown rows (now rows 30-34).
$ws.Range("A30:A34").EntireRow.Hidden = $true

# 9. Update sheet view: drop frozen/scrolled topLeftCell, move selection to J26
[void]$ws.Range("J26").Select()

Write-Host "done"
